$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10773.9
$ws.Range("I9").Value = 15219.571
$ws.Range("J9").Value = 400.66666
$ws.Range("K9").Value = 15219.571
$ws.Range("L9").Value = 400.66666
$ws.Range("M9").Value = -15050.571
$ws.Range("N9").Value = -738.66666

$ws.Range("H86").Value = 6282
$ws.Range("I86").Value = 8424.25
$ws.Range("K86").Value = 8424.25
$ws.Range("M86").Value = -7301.25

$ws.Range("H89").Value = 6282
$ws.Range("I89").Value = 8424.25
$ws.Range("K89").Value = 42121.25
$ws.Range("M89").Value = -36505.25

$ws.Range("H98").Value = 1139.0588
$ws.Range("I98").Value = 870.5454999999999
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 870.5454999999999
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 627.4545000000001
$ws.Range("N98").Value = -12996

$ws.Range("H112").Value = 97487.95
$ws.Range("I112").Value = 144199.72
$ws.Range("J112").Value = 74132.07000000001
$ws.Range("K112").Value = 432599.16
$ws.Range("L112").Value = 222396.21
$ws.Range("M112").Value = -431491.16
$ws.Range("N112").Value = -224612.21

$ws.Range("H122").Value = 1139.0588
$ws.Range("I122").Value = 870.5454999999999
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 2611.6365
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -161.6364999999996
$ws.Range("N122").Value = -34900

$ws.Range("H138").Value = 1722.7931
$ws.Range("I138").Value = 1170.6562
$ws.Range("K138").Value = 3511.9686
$ws.Range("M138").Value = 1628.0314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 206
$ws.Range("J3").Value = 206
$ws.Range("L3").Value = 206
$ws.Range("N3").Value = -436

$ws.Range("H5").Value = 721.6
$ws.Range("I5").Value = 689.625
$ws.Range("J5").Value = 849.5
$ws.Range("K5").Value = 689.625
$ws.Range("L5").Value = 849.5
$ws.Range("M5").Value = -577.625
$ws.Range("N5").Value = -1073.5

$ws.Range("H45").Value = 5308
$ws.Range("I45").Value = 5246.273
$ws.Range("J45").Value = 5477.75
$ws.Range("K45").Value = 5246.273
$ws.Range("L45").Value = 5477.75
$ws.Range("M45").Value = -4869.273
$ws.Range("N45").Value = -6231.75

$ws.Range("H74").Value = 27029414
$ws.Range("I74").Value = 31252340
$ws.Range("J74").Value = 2698.2
$ws.Range("K74").Value = 31252340
$ws.Range("L74").Value = 2698.2
$ws.Range("M74").Value = -31251466
$ws.Range("N74").Value = -4446.2

$ws.Range("H77").Value = 27029414
$ws.Range("I77").Value = 31252340
$ws.Range("J77").Value = 2698.2
$ws.Range("K77").Value = 156261700
$ws.Range("L77").Value = 13491
$ws.Range("M77").Value = -156257332
$ws.Range("N77").Value = -22227

$ws.Range("H110").Value = 127086.625
$ws.Range("I110").Value = 167782.17
$ws.Range("K110").Value = 167782.17
$ws.Range("M110").Value = -165737.17

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 721.6
$ws.Range("I4").Value = 689.625
$ws.Range("J4").Value = 849.5
$ws.Range("K4").Value = 689.625
$ws.Range("L4").Value = 849.5
$ws.Range("M4").Value = -574.625
$ws.Range("N4").Value = -1079.5

$ws.Range("H5").Value = 836
$ws.Range("I5").Value = 1469.5
$ws.Range("J5").Value = 202.5
$ws.Range("K5").Value = 1469.5
$ws.Range("L5").Value = 202.5
$ws.Range("M5").Value = -1356.5
$ws.Range("N5").Value = -428.5

$ws.Range("H86").Value = 3340.1936
$ws.Range("I86").Value = 3267.3044
$ws.Range("J86").Value = 3549.75
$ws.Range("K86").Value = 3267.3044
$ws.Range("L86").Value = 3549.75
$ws.Range("M86").Value = -2144.3044
$ws.Range("N86").Value = -5795.75

$ws.Range("H89").Value = 3340.1936
$ws.Range("I89").Value = 3267.3044
$ws.Range("J89").Value = 3549.75
$ws.Range("K89").Value = 16336.522
$ws.Range("L89").Value = 17748.75
$ws.Range("M89").Value = -10720.522
$ws.Range("N89").Value = -28980.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2331
$ws.Range("I16").Value = 2179.7273
$ws.Range("K16").Value = 2179.7273
$ws.Range("M16").Value = -1892.7273

$ws.Range("H68").Value = 98998
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 98998
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 98998
$ws.Range("N68").Value = -100496
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 98998
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 98998
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 296994
$ws.Range("N71").Value = -304482
$ws.Range("M71").ClearContents()

$ws.Range("H95").Value = 17262
$ws.Range("J95").Value = 17262
$ws.Range("L95").Value = 17262
$ws.Range("N95").Value = -22754

$ws.Range("H113").Value = 2331
$ws.Range("I113").Value = 2179.7273
$ws.Range("K113").Value = 2179.7273
$ws.Range("M113").Value = -9.727300000000014

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1094.7142
$ws.Range("I9").Value = 1229
$ws.Range("J9").Value = 994
$ws.Range("K9").Value = 3687
$ws.Range("L9").Value = 2982
$ws.Range("M9").Value = -3463
$ws.Range("N9").Value = -3430

$ws.Range("H114").Value = 112243.78
$ws.Range("I114").Value = 125774.25
$ws.Range("J114").Value = 4000
$ws.Range("K114").Value = 377322.75
$ws.Range("L114").Value = 12000
$ws.Range("M114").Value = -374068.75
$ws.Range("N114").Value = -18508

$ws.Range("H120").Value = 28331.334
$ws.Range("I120").Value = 27497.5
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 82492.5
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -77654.5
$ws.Range("N120").Value = -99673

$ws.Range("H132").Value = 1722.9
$ws.Range("J132").Value = 1891.96
$ws.Range("L132").Value = 17027.64
$ws.Range("N132").Value = -22087.64

$ws.Range("H138").Value = 3154.2
$ws.Range("I138").Value = 3154.2
$ws.Range("K138").Value = 9462.599999999999
$ws.Range("M138").Value = -4322.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 99999
$ws.Range("J4").Value = 99999
$ws.Range("L4").Value = 99999
$ws.Range("N4").Value = -100223

$ws.Range("H97").Value = 1759.4445
$ws.Range("I97").Value = 1725.6428
$ws.Range("J97").Value = 1877.75
$ws.Range("K97").Value = 1725.6428
$ws.Range("L97").Value = 1877.75
$ws.Range("M97").Value = -1229.6428
$ws.Range("N97").Value = -2869.75

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H113").Value = 22079.691
$ws.Range("I113").Value = 27628.049
$ws.Range("K113").Value = 27628.049
$ws.Range("M113").Value = -25458.049

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3490.6365
$ws.Range("I46").Value = 2660
$ws.Range("J46").Value = 4182.8335
$ws.Range("K46").Value = 2660
$ws.Range("L46").Value = 4182.8335
$ws.Range("M46").Value = -2472
$ws.Range("N46").Value = -4558.8335

$ws.Range("H68").Value = 3345.1667
$ws.Range("J68").Value = 3083
$ws.Range("L68").Value = 3083
$ws.Range("N68").Value = -4581

$ws.Range("H71").Value = 3345.1667
$ws.Range("J71").Value = 3083
$ws.Range("L71").Value = 15415
$ws.Range("N71").Value = -22903

$ws.Range("H82").Value = 1377.4
$ws.Range("I82").Value = 1377.4
$ws.Range("K82").Value = 1377.4
$ws.Range("M82").Value = -1016.4

$ws.Range("H85").Value = 1377.4
$ws.Range("I85").Value = 1377.4
$ws.Range("K85").Value = 1377.4
$ws.Range("M85").Value = -129.4000000000001

$ws.Range("H93").Value = 2033.1111
$ws.Range("I93").Value = 1256.8572
$ws.Range("J93").Value = 4750
$ws.Range("K93").Value = 1256.8572
$ws.Range("L93").Value = 4750
$ws.Range("M93").Value = -8.857199999999921
$ws.Range("N93").Value = -7246

$ws.Range("H122").Value = 8524.916999999999
$ws.Range("I122").Value = 5922.222
$ws.Range("J122").Value = 16333
$ws.Range("K122").Value = 17766.666
$ws.Range("L122").Value = 48999
$ws.Range("M122").Value = -15316.666
$ws.Range("N122").Value = -53899

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1420.4762
$ws.Range("I113").Value = 1460.3334
$ws.Range("J113").Value = 1320.8334
$ws.Range("K113").Value = 4381.0002
$ws.Range("L113").Value = 3962.5002
$ws.Range("M113").Value = -2211.0002
$ws.Range("N113").Value = -8302.5002

$ws.Range("H132").Value = 16675265
$ws.Range("I132").Value = 25003822
$ws.Range("K132").Value = 75011466
$ws.Range("M132").Value = -75008936
